# Updated cryptos list on Sun Aug 13 23:41:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep being stored as text (matches source
# data format, e.g. "29.311.73", "0.9990", "  -0.44%  ") instead of being
# auto-coerced into numbers when the new values are assigned below.
$ws.Range("D2:E51").NumberFormat = "@"

# Price (column D) and Volume(1h) (column E) updates per row
$ws.Range("D2").Value = "29.311.73"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.841.35"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").Value = "0.9983"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "240.24"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "0.6265"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "0.07484"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("D9").Value = "0.2893"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  -2.20%  "

$ws.Range("D11").Value = "0.07715"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").Value = "1.839.39"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").Value = "0.6776"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "0.00001027"
$ws.Range("E15").Value = "  -4.32%  "

$ws.Range("D16").Value = "82.12"
$ws.Range("E16").Value = "  -1.68%  "

$ws.Range("D17").Value = "2.100.34"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").Value = "6.101"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").Value = "29.313.49"

$ws.Range("D20").Value = "228.60"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "7.374"
$ws.Range("E23").Value = "  -1.00%  "

$ws.Range("D24").Value = "0.9990"
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").Value = "158.28"
$ws.Range("E25").Value = "  +0.26%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "8.373"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("D28").Value = "17.55"
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").Value = "1.395"
$ws.Range("E29").Value = "  +2.00%  "

$ws.Range("D30").Value = "1.472"
$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("D31").Value = "0.05706"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").Value = "4.099"
$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("D33").Value = "4.032"
$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("E34").Value = "  -1.43%  "

$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("D36").Value = "0.6929"
$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("D37").Value = "2.584"
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").Value = "2.820"
$ws.Range("E38").Value = "  +2.33%  "

$ws.Range("D39").Value = "1.244.32"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("D40").Value = "0.01811"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("D42").Value = "0.9037"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").Value = "0.9985"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").Value = "2.000.87"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("D45").Value = "101.39"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "65.76"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").Value = "7.069"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("D48").Value = "0.1168"
$ws.Range("E48").Value = "  +1.20%  "

# Rows 49-51: BabyDogeCoin dropped from list, remaining coins shift up,
# RenderToken added as the new last row
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.968"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.3938"
$ws.Range("E50").Value = "  -2.09%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.657"
$ws.Range("E51").Value = "  -1.13%  "
